# StudentTemplate.xlsx: import-student template now supports adding one or
# more "Groepnaam" (group name) values, so a new header column E is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell - this also grows the shared-string table and the sheet's
# used-range dimension to A1:E1 automatically.
$ws.Range("E1").Value = "Groepnaam"

# Give the new column a sensible, explicit width (author manually resized it).
$ws.Columns("E").ColumnWidth = 13

# Move/collapse the selection to E2, as left by the author after editing E1.
$ws.Range("E2").Select()
